$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "71.158.71"
$ws.Cells.Item(2, 5).Value = "  +0.30%  "
$ws.Cells.Item(3, 4).Value = "3.870.52"
$ws.Cells.Item(3, 5).Value = "  +1.64%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "698.04"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +0.31%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "173.90"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +0.40%  "
$ws.Cells.Item(7, 4).Value = "3.868.10"
$ws.Cells.Item(7, 5).Value = "  +1.67%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$ws.Cells.Item(9, 5).Value = "  +0.26%  "
$ws.Cells.Item(10, 5).Value = "  -0.07%  "
$ws.Cells.Item(11, 5).Value = "  -4.98%  "
$ws.Cells.Item(12, 5).Value = "  -0.16%  "
$ws.Cells.Item(13, 5).Value = "  +3.40%  "
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "36.49"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +0.72%  "
$ws.Cells.Item(15, 4).Value = "4.525.02"
$ws.Cells.Item(15, 5).Value = "  +1.60%  "
$ws.Cells.Item(16, 4).Value = "3.873.33"
$ws.Cells.Item(16, 5).Value = "  +1.65%  "
$ws.Cells.Item(17, 4).Value = "71.248.90"
$ws.Cells.Item(17, 5).Value = "  +0.39%  "
$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.25"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +0.69%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.70"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -0.41%  "
$ws.Cells.Item(20, 5).Value = "  -0.40%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.17"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -1.43%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "500.26"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +4.47%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.725"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +1.40%  "
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "84.94"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +1.53%  "
$ws.Cells.Item(25, 2).Value = "PEPE"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000148"
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.87"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +5.43%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.28"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -0.82%  "
$ws.Cells.Item(28, 5).Value = "  +1.04%  "
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.20"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +3.04%  "
$ws.Cells.Item(30, 5).Value = "  -0.03%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +1.60%  "
$ws.Cells.Item(32, 5).Value = "  -1.30%  "
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.75"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +0.55%  "
$ws.Cells.Item(34, 5).Value = "  +2.30%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.26"
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +0.78%  "
$ws.Cells.Item(36, 4).Value = "3.822.78"
$ws.Cells.Item(36, 5).Value = "  +1.68%  "
$ws.Cells.Item(37, 5).Value = "  +1.13%  "
$ws.Cells.Item(38, 5).Value = "  +2.63%  "
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.41"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +10.05%  "
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.43"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -2.88%  "
$ws.Cells.Item(41, 5).Value = "  +8.37%  "
$ws.Cells.Item(42, 5).Value = "  +1.15%  "
$ws.Cells.Item(43, 5).Value = "  +0.02%  "
$ws.Cells.Item(44, 5).Value = "  -0.02%  "
$ws.Cells.Item(45, 2).Value = "Monero"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "163.82"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +2.09%  "
$ws.Cells.Item(46, 2).Value = "FLOKI"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.000311"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -6.80%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "48.88"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -1.24%  "
$ws.Cells.Item(48, 5).Value = "  +1.75%  "
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "417.79"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +4.79%  "
$ws.Cells.Item(50, 5).Value = "  -2.62%  "
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "43.73"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -2.80%  "
